# Apple.xlsx edit: fill in the "Relevância" (column B) labels that were
# typed in for more rows of the "Treinamento" and "Teste" sheets, and
# update each sheet's saved view/selection state to match where the
# author was last working.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Treinamento"
$ws2 = $wb.Worksheets.Item(2)   # "Teste"

# ---------------------------------------------------------------------
# Sheet "Treinamento": new B172:B180 relevance labels
# ---------------------------------------------------------------------
$treinoVals = @(1,0,1,0,1,0,0,0,1)
$treinoStartRow = 172
for ($i = 0; $i -lt $treinoVals.Length; $i++) {
    $ws1.Cells.Item($treinoStartRow + $i, 2).Value = $treinoVals[$i]
}

# ---------------------------------------------------------------------
# Sheet "Teste": new B150:B201 relevance labels
# ---------------------------------------------------------------------
$testeVals = @(1,1,0,1,0,1,0,1,1,1,1,1,0,0,1,1,1,0,0,1,1,0,0,1,0,0,1,0,1,1,1,1,1,0,0,0,1,1,1,1,0,1,1,0,0,0,0,1,0,1,0,1)
$testeStartRow = 150
for ($i = 0; $i -lt $testeVals.Length; $i++) {
    $ws2.Cells.Item($testeStartRow + $i, 2).Value = $testeVals[$i]
}

# ---------------------------------------------------------------------
# View/selection state: "Teste" was scrolled/selected while editing,
# then the author came back to "Treinamento" leaving it as the active
# (saved) tab with its own scroll position/selection.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B202:B207").Select()

$ws1.Activate()
$ws1.Range("B181").Select()
